# Insert a new data row at row 263 (weekly price report gained a new entry),
# pushing the existing rows 263-368 down to 264-369.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(263).Insert()

# Populate the newly inserted row 263 with the new record.
$ws.Range("A263").Value = 3
$ws.Range("B263").Value = "Femacal de La Calera"
$ws.Range("C263").Value = "Coquimbo"
$ws.Range("D263").Value2 = 44755
$ws.Range("E263").Value = 5
$ws.Range("F263").Value = 100112012
$ws.Range("G263").Value = "Espinaca"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 120
$ws.Range("K263").Value = 4000
$ws.Range("L263").Value = 4000
$ws.Range("M263").Value = 4000
$ws.Range("N263").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O263").Value = "Provincia de Quillota"
$ws.Range("P263").Value = 1333
$ws.Range("Q263").Value = 3
$ws.Range("R263").Value = "Hortaliza"
